$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "string" table is sorted alphabetically by the "keys" column (A).
# Two new translation keys - "year_overview" and "yearly_average" - are
# inserted between the existing "year" (row 106) and "yes" (row 107) rows,
# so insert two blank rows at 107 (this pushes "yes" down to row 109).
$ws.Range("A107:A108").EntireRow.Insert() | Out-Null

# Populate the two new rows. Values are written in this particular order
# (not simply left-to-right, row-by-row) so that the workbook's shared
# string table ends up listing the new unique strings in the same order
# they were originally authored.
$ws.Range("A107").Value = "year_overview"
$ws.Range("D107").Value = "Jahresübersicht"
$ws.Range("E107").Value = "Joresiwwersiicht"
$ws.Range("C107").Value = "Aperçu de l'année"

$ws.Range("A108").Value = "yearly_average"
$ws.Range("B108").Value = "Yearly average"
$ws.Range("B107").Value = "Year overview"
$ws.Range("C108").Value = "Moyenne annuelle"
$ws.Range("D108").Value = "Jahresdurchschnitt"
$ws.Range("E108").Value = "Joresmoyenne"

# Grow the "string" table so it covers the two extra rows (A1:E107 -> A1:E109).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E109")) | Out-Null

# The ExternalData_1 defined name tracked the table's first two columns -
# extend it to match the new last row.
$wb.Names.Item(1).RefersTo = "=string!`$A`$1:`$B`$109"

# Reflect the author's final cursor position in the saved view state.
$ws.Range("B93").Select() | Out-Null
